$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking value updates
$ws.Range('D2').Value = '43.837.35'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '2.315.81'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('E5').Value = '  +20.13%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('E10').Value = '  +7.93%  '
$ws.Range('E11').Value = '  +1.75%  '
$ws.Range('E12').Value = '  +16.52%  '
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('E14').Value = '  +3.38%  '
$ws.Range('D15').Value = '2.661.48'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '2.307.71'
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '43.829.92'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('E19').Value = '  +3.46%  '
$ws.Range('E20').Value = '  +8.86%  '
$ws.Range('E21').Value = '  +1.14%  '
$ws.Range('E22').Value = '  +6.15%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  +17.69%  '
$ws.Range('E25').Value = '  +6.98%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('E28').Value = '  +9.46%  '
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  +5.46%  '
$ws.Range('E34').Value = '  +4.38%  '
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('E36').Value = '  +6.53%  '
$ws.Range('E37').Value = '  +3.32%  '
$ws.Range('E38').Value = '  +21.88%  '
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('E40').Value = '  +3.59%  '
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('E42').Value = '  +13.51%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('E43').Value = '  +7.91%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('B45').Value = 'THORChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('E45').Value = '  +14.70%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('E49').Value = '  +11.51%  '
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('E51').Value = '  +3.84%  '

# Numeric-looking values that must remain stored as text (match source formatting)
# Force text storage via NumberFormat "@" then write value, then clear the temporary
# number-format override afterwards so the cell style matches the original (no explicit style).
$forceTextCells = @('D5', 'D6', 'D9', 'D10', 'D11', 'D12', 'D14', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D36', 'D39', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range('D5').Value = '114.63'
$ws.Range('D6').Value = '270.60'
$ws.Range('D9').Value = '0.624'
$ws.Range('D10').Value = '47.75'
$ws.Range('D11').Value = '0.0946'
$ws.Range('D12').Value = '8.91'
$ws.Range('D14').Value = '15.69'
$ws.Range('D19').Value = '0.0000111'
$ws.Range('D20').Value = '6.71'
$ws.Range('D21').Value = '72.83'
$ws.Range('D22').Value = '2.50'
$ws.Range('D23').Value = '234.45'
$ws.Range('D24').Value = '2.94'
$ws.Range('D25').Value = '9.54'
$ws.Range('D27').Value = '11.49'
$ws.Range('D28').Value = '42.60'
$ws.Range('D31').Value = '177.85'
$ws.Range('D32').Value = '21.99'
$ws.Range('D33').Value = '0.0933'
$ws.Range('D34').Value = '5.59'
$ws.Range('D36').Value = '4.73'
$ws.Range('D39').Value = '0.0358'
$ws.Range('D42').Value = '70.78'
$ws.Range('D43').Value = '12.79'
$ws.Range('D44').Value = '1.00'
$ws.Range('D45').Value = '5.97'
$ws.Range('D46').Value = '1.40'
$ws.Range('D47').Value = '8.84'
$ws.Range('D49').Value = '0.470'
$ws.Range('D50').Value = '100.57'
$ws.Range('D51').Value = '1.24'
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).ClearFormats()
}
